$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 2507.2856
$ws.Range("I80").Value = 220
$ws.Range("J80").Value = 3778
$ws.Range("K80").Value = 660
$ws.Range("L80").Value = 11334
$ws.Range("M80").Value = 338
$ws.Range("N80").Value = -13330
$ws.Range("H83").Value = 2507.2856
$ws.Range("I83").Value = 220
$ws.Range("J83").Value = 3778
$ws.Range("K83").Value = 1980
$ws.Range("L83").Value = 34002
$ws.Range("M83").Value = 3012
$ws.Range("N83").Value = -43986
$ws.Range("H96").Value = 930.7692
$ws.Range("I96").Value = 589.3333
$ws.Range("K96").Value = 1767.9999
$ws.Range("M96").Value = -394.9999
$ws.Range("H106").Value = 9370.166999999999
$ws.Range("I106").Value = 1731.25
$ws.Range("K106").Value = 1731.25
$ws.Range("M106").Value = -1100.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 5973.75
$ws.Range("I122").Value = 895
$ws.Range("J122").Value = 7666.6665
$ws.Range("K122").Value = 2685
$ws.Range("L122").Value = 22999.9995
$ws.Range("M122").Value = -235
$ws.Range("N122").Value = -27899.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 916.3333
$ws.Range("I80").Value = 922
$ws.Range("J80").Value = 914.44446
$ws.Range("K80").Value = 922
$ws.Range("L80").Value = 914.44446
$ws.Range("M80").Value = 76
$ws.Range("N80").Value = -2910.44446
$ws.Range("H83").Value = 916.3333
$ws.Range("I83").Value = 922
$ws.Range("J83").Value = 914.44446
$ws.Range("K83").Value = 4610
$ws.Range("L83").Value = 4572.2223
$ws.Range("M83").Value = 382
$ws.Range("N83").Value = -14556.2223
$ws.Range("H94").Value = 450.26086
$ws.Range("I94").Value = 571.8461
$ws.Range("J94").Value = 292.2
$ws.Range("K94").Value = 571.8461
$ws.Range("L94").Value = 292.2
$ws.Range("M94").Value = -120.8461
$ws.Range("N94").Value = -1194.2
$ws.Range("H99").Value = 1089.8948
$ws.Range("I99").Value = 891.63635
$ws.Range("K99").Value = 891.63635
$ws.Range("M99").Value = 606.36365
$ws.Range("H105").Value = 7118.5454
$ws.Range("J105").Value = 10963.462
$ws.Range("L105").Value = 10963.462
$ws.Range("N105").Value = -14457.462
$ws.Range("H107").Value = 1900
$ws.Range("I107").Value = 1900
$ws.Range("K107").Value = 1900
$ws.Range("M107").Value = 20
$ws.Range("H132").Value = 69999
$ws.Range("J132").Value = 69999
$ws.Range("L132").Value = 69999
$ws.Range("N132").Value = -80119

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2013
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("H107").Value = 909.1852
$ws.Range("I107").Value = 803.04346
$ws.Range("K107").Value = 803.04346
$ws.Range("M107").Value = 1116.95654
$ws.Range("H113").Value = 2013
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 5084
$ws.Range("I68").Value = 4175
$ws.Range("J68").Value = 5311.25
$ws.Range("K68").Value = 12525
$ws.Range("L68").Value = 15933.75
$ws.Range("M68").Value = -11714
$ws.Range("N68").Value = -17555.75
$ws.Range("H71").Value = 5084
$ws.Range("I71").Value = 4175
$ws.Range("J71").Value = 5311.25
$ws.Range("K71").Value = 37575
$ws.Range("L71").Value = 47801.25
$ws.Range("M71").Value = -33519
$ws.Range("N71").Value = -55913.25
$ws.Range("H132").Value = 4549.9414
$ws.Range("J132").Value = 6281.25
$ws.Range("L132").Value = 56531.25
$ws.Range("N132").Value = -61591.25
$ws.Range("H137").Value = 5364.615
$ws.Range("I137").Value = 1663.1428
$ws.Range("K137").Value = 4989.428400000001
$ws.Range("M137").Value = 110.5715999999993

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 6712.8
$ws.Range("I122").Value = 2852
$ws.Range("J122").Value = 12504
$ws.Range("K122").Value = 8556
$ws.Range("L122").Value = 37512
$ws.Range("M122").Value = -6106
$ws.Range("N122").Value = -42412
$ws.Range("H132").Value = 45915.73
$ws.Range("I132").Value = 64963.117
$ws.Range("K132").Value = 194889.351
$ws.Range("M132").Value = -192359.351

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 13325.75
$ws.Range("I7").Value = 3075.8
$ws.Range("J7").Value = 30409
$ws.Range("K7").Value = 3075.8
$ws.Range("L7").Value = 30409
$ws.Range("M7").Value = -2963.8
$ws.Range("N7").Value = -30633
$ws.Range("H40").Value = 6868.5884
$ws.Range("I40").Value = 8321.5
$ws.Range("K40").Value = 8321.5
$ws.Range("M40").Value = -8185.5
$ws.Range("H61").Value = 6821.2144
$ws.Range("J61").Value = 12519
$ws.Range("L61").Value = 12519
$ws.Range("N61").Value = -12923
$ws.Range("H113").Value = 6821.2144
$ws.Range("J113").Value = 12519
$ws.Range("L113").Value = 12519
$ws.Range("N113").Value = -16859
$ws.Range("H126").Value = 13325.75
$ws.Range("I126").Value = 3075.8
$ws.Range("J126").Value = 30409
$ws.Range("K126").Value = 9227.400000000001
$ws.Range("L126").Value = 91227
$ws.Range("M126").Value = -6757.400000000001
$ws.Range("N126").Value = -96167

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 1903199.8
$ws.Range("I29").Value = 1553999.5
$ws.Range("J29").Value = 2136000
$ws.Range("K29").Value = 1553999.5
$ws.Range("L29").Value = 2136000
$ws.Range("M29").Value = -1553709.5
$ws.Range("N29").Value = -2136580
$ws.Range("H81").Value = 4106.222
$ws.Range("I81").Value = 2314.5454
$ws.Range("K81").Value = 4629.0908
$ws.Range("M81").Value = -3568.0908
$ws.Range("H84").Value = 4106.222
$ws.Range("I84").Value = 2314.5454
$ws.Range("K84").Value = 23145.454
$ws.Range("M84").Value = -17841.454
$ws.Range("H113").Value = 779.6667
$ws.Range("I113").Value = 735.7273
$ws.Range("J113").Value = 900.5
$ws.Range("K113").Value = 2207.1819
$ws.Range("L113").Value = 2701.5
$ws.Range("M113").Value = -37.18190000000004
$ws.Range("N113").Value = -7041.5
$ws.Range("H132").Value = 3352.8057
$ws.Range("I132").Value = 3184.4062
$ws.Range("K132").Value = 9553.2186
$ws.Range("M132").Value = -7023.2186
$ws.Range("H141").Value = 97554
$ws.Range("J141").Value = 97554
$ws.Range("L141").Value = 97554
$ws.Range("N141").Value = -107914
